$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nt_POT-GPD-Equivalent")

# --- Update existing rows 5-22 (column B text content changes) ---
$ws.Cells.Item(5, 2).Value = '[3] "      estimate Std. Error"                                   '
$ws.Cells.Item(6, 2).Value = '[4] "dummy        1         NA"                                   '
$ws.Cells.Item(7, 2).Value = '[5] "Fixed parameters:"                                           '
$ws.Cells.Item(8, 2).Value = '[6] "            value"                                           '
$ws.Cells.Item(9, 2).Value = '[7] "scale    23.39786"                                           '
$ws.Cells.Item(10, 2).Value = '[8] "location 42.00000"                                           '
$ws.Cells.Item(11, 2).Value = '[9] "shape     0.00001"                                           '
$ws.Range("B12").ClearContents()
$ws.Cells.Item(13, 2).Value = '$`GoodnessOfFit_fitdistrplus::gofstat`'
$ws.Cells.Item(14, 2).Value = ' [1] "Goodness-of-fit statistics"                '
$ws.Cells.Item(15, 2).Value = ' [2] "                             1-mle-mygpd"  '
$ws.Cells.Item(16, 2).Value = ' [3] "Kolmogorov-Smirnov statistic   0.0889465"  '
$ws.Cells.Item(17, 2).Value = ' [4] "Cramer-von Mises statistic     0.2092053"  '
$ws.Cells.Item(18, 2).Value = ' [5] "Anderson-Darling statistic     1.6838894"  '
$ws.Cells.Item(19, 2).Value = ' [6] ""                                          '
$ws.Cells.Item(20, 2).Value = ' [7] "Goodness-of-fit criteria"                  '
$ws.Cells.Item(21, 2).Value = ' [8] "                               1-mle-mygpd"'
$ws.Cells.Item(22, 2).Value = ' [9] "Akaike''s Information Criterion    1505.242"'

# --- Column A for new rows 23-33: copy sequential-number strings from a sister sheet ---
$wsSrc = $wb.Worksheets.Item("nt_distLquantile_parameters")
$wsSrc.Range("A23:A33").Copy($ws.Range("A23"))

# --- Column B for new rows 23-33 ---
$ws.Cells.Item(23, 2).Value = '[10] "Bayesian Information Criterion    1508.441"'
$ws.Cells.Item(25, 2).Value = '$`KolmogorovSmirnovTest_stats::ks.test`'
$ws.Cells.Item(26, 2).Value = '[1] ""                                    "\tOne-sample Kolmogorov-Smirnov test"'
$ws.Cells.Item(27, 2).Value = '[3] ""                                    "data:  imp.vals$nt.series"          '
$ws.Cells.Item(28, 2).Value = '[5] "D = 0.088947, p-value = 0.1141"      "alternative hypothesis: two-sided"  '
$ws.Cells.Item(29, 2).Value = '[7] ""                                   '
$ws.Cells.Item(31, 2).Value = '$RMSE'
$ws.Cells.Item(32, 2).Value = '[1] "[1] 0.02322051"'
